# Fruta / hortaliza, semanal
# A new daily price record was inserted into the "Kiwi" price log at row 739,
# pushing every existing row from 739 downward by one (old row 845 becomes
# the new row 846). This mirrors the author's weekly data refresh where a
# new day's quote is prepended to this product's block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh blank row at position 739; Excel shifts 739:845 -> 740:846.
$ws.Rows("739:739").Insert()

# Populate the newly inserted row with the new record's values.
$ws.Cells.Item(739, 1).Value = 10
$ws.Cells.Item(739, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(739, 3).Value = "La Araucanía"
$ws.Cells.Item(739, 4).Value = 45154
$ws.Cells.Item(739, 5).Value = 9
$ws.Cells.Item(739, 6).Value = "Fruta"
$ws.Cells.Item(739, 7).Value = 100101
$ws.Cells.Item(739, 8).Value = "Berries"
$ws.Cells.Item(739, 9).Value = 100101007
$ws.Cells.Item(739, 10).Value = "Kiwi"
$ws.Cells.Item(739, 11).Value = "Hayward"
$ws.Cells.Item(739, 12).Value = "Especial"
$ws.Cells.Item(739, 13).Value = 110
$ws.Cells.Item(739, 14).Value = 16000
$ws.Cells.Item(739, 15).Value = 16000
$ws.Cells.Item(739, 16).Value = 16000
$ws.Cells.Item(739, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(739, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(739, 19).Value = 1600
$ws.Cells.Item(739, 20).Value = 10
